$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "hello"
Write-Host "A1: $($ws.Range('A1').Value)"
